# feat: add 2022-Q1 data
#
# - Inserts a new worksheet "2022-Q1" (fund-holdings detail, same layout as
#   "2021-Q4") between the existing "2021-Q4" and "总计" sheets.
# - Updates the "总计" (totals) sheet with a new row summarising 2022-Q1
#   (8 holdings, 9.46 亿元) ahead of the existing 2021-Q4 row.

function Set-TextCell {
    param($Cell, $Text)
    # Force a numeric-looking string ("004997", "148.04", ...) to be stored
    # as text instead of being auto-coerced to a number, while leaving the
    # cell on the default ("Normal") style afterwards - matches the source
    # cells (t="inlineStr", no explicit s="..").
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# --- 1. Insert the new "2022-Q1" sheet right before "总计" ----------------
# Duplicate "总计" itself (rather than Worksheets.Add(), which creates a
# bare sheet missing the workbook's usual <sheetPr> block) so the new sheet
# starts out with the same sheetPr/page-setup boilerplate as its siblings;
# its cell content gets fully overwritten below anyway.
$total.Copy($total, $null)
$q1 = $wb.Worksheets.Item(2)
$q1.Name = "2022-Q1"

# NOTE: after Copy()/rename, the old `$total` variable tracks the *slot*
# (now occupied by the new "2022-Q1" copy) rather than following the
# original "总计" worksheet to its new position - re-resolve it by name.
$total = $wb.Worksheets.Item("总计")

# Pull over the header-row style (bold/centered/bordered = style index 2 in
# the source) and the column-A index style from the existing 2021-Q4 sheet.
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q4.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)

# Clear any leftover "总计"-sized values outside the new A1:H9 footprint
# (there are none beyond D2 in the source, but be explicit/safe).
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund holdings detail rows (index, code, name, scale, total position,
# position share, held value 亿元, position rank)
$fundRows = @(
    @(0, "004997", "广发高端制造股票A",         "148.04", "94.19", "4.31", "6.3805", 10),
    @(1, "011479", "广发诚享混合A",              "44.44",  "93.13", "4.48", "1.9909", 10),
    @(2, "010160", "广发高端制造股票C",          "8.51",   "94.19", "4.31", "0.3668", 10),
    @(3, "001468", "广发改革先锋灵活配置混合",   "10.51",  "76.65", "3.47", "0.3647", 3),
    @(4, "011480", "广发诚享混合C",              "4.52",   "93.13", "4.48", "0.2025", 10),
    @(5, "900029", "中信证券量化优选股票A",      "3.91",   "90.12", "2.09", "0.0817", 9),
    @(6, "900030", "中信证券量化优选股票C",      "1.94",   "90.12", "2.09", "0.0405", 9),
    @(7, "162717", "广发再融资主题灵活配置混合", "0.65",   "79.18", "4.92", "0.0320", 5)
)

$r = 2
foreach ($row in $fundRows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    Set-TextCell $q1.Cells.Item($r, 2) $row[1]
    Set-TextCell $q1.Cells.Item($r, 3) $row[2]
    Set-TextCell $q1.Cells.Item($r, 4) $row[3]
    Set-TextCell $q1.Cells.Item($r, 5) $row[4]
    Set-TextCell $q1.Cells.Item($r, 6) $row[5]
    Set-TextCell $q1.Cells.Item($r, 7) $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# --- 2. Add the 2022-Q1 summary row to "总计" ------------------------------
# Push the existing 2021-Q4 summary row down to row 3 (copy its current
# values/format as-is) and overwrite row 2 with the new 2022-Q1 summary
# (most recent quarter first). Avoid Rows.Insert() here - it copies the
# bold/bordered header-row formatting down onto the blank row, which is
# not what the plain default-styled data rows need.
$total.Range("A2:D2").Copy()
$total.Range("A3:D3").PasteSpecial(-4122)
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 10.05
$total.Range("A3").Value = 1

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 9.46
